$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally lists metrics for 8 models (rows 2-8):
# Logistic Regression, Lasso, Support Vector Classifier, CART,
# Random Forest, LightGBM, XGBoost.
# Only "Logistic Regression" and "LightGBM" remain relevant, so:
#   - row 2 (Logistic Regression) is kept as-is
#   - row 3 is overwritten with the LightGBM figures (previously row 7)
#   - rows 4-8 (Lasso's old data plus SVC/CART/RandomForest/XGBoost rows) are removed

# Capture the LightGBM row (row 7) values before we start shifting rows around.
$modelName = $ws.Cells.Item(7, 1).Value()
$accuracy = $ws.Cells.Item(7, 2).Value()
$precisionWeighted = $ws.Cells.Item(7, 3).Value()
$recallWeighted = $ws.Cells.Item(7, 4).Value()
$f1Weighted = $ws.Cells.Item(7, 5).Value()

# Overwrite row 3 (currently "Lasso") with the LightGBM data.
$ws.Cells.Item(3, 1).Value = $modelName
$ws.Cells.Item(3, 2).Value = $accuracy
$ws.Cells.Item(3, 3).Value = $precisionWeighted
$ws.Cells.Item(3, 4).Value = $recallWeighted
$ws.Cells.Item(3, 5).Value = $f1Weighted

# Remove the now-redundant rows 4 through 8 entirely.
$ws.Range("A4:E8").EntireRow.Delete()
